$wb = $excel.ActiveWorkbook

# --- AppOvfProperties: insert a "DataCentreName" column at the front ---
$wsApp = $wb.Worksheets.Item("AppOvfProperties")
$wsApp.Columns.Item(1).Insert()
$wsApp.Range("A1").Value = "DataCentreName"
$wsApp.Range("A2:A5").Value = "FakeCorp"
$wsApp.Columns.Item(1).AutoFit()
[void]$wsApp.Range("A1:A5").Select()

# --- VmOvfProperties: insert a "DataCentreName" column at the front ---
$wsVm = $wb.Worksheets.Item("VmOvfProperties")
$wsVm.Columns.Item(1).Insert()
$wsVm.Range("A1").Value = "DataCentreName"
$wsVm.Range("A2:A5").Value = "FakeCorp"

# --- Imports: selection moved to full column C ---
$wsImports = $wb.Worksheets.Item("Imports")
[void]$wsImports.Range("C:C").Select()

# --- AppExports: selection moved to A2 ---
$wsExports = $wb.Worksheets.Item("AppExports")
[void]$wsExports.Range("A2").Select()

# VmOvfProperties ends up the active sheet, selection at C12
$wsVm.Activate()
[void]$wsVm.Range("C12").Select()
